$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3, column B currently holds "JOHN DEERE" (all caps) - update to proper case "John Deere"
$ws.Range("B3").Value = "John Deere"
